# Applies the "Optuna Attempt (go back with original)" edit:
# updates MyForecast / Inventory Coverage / Stockout Risk / Seasonality Index
# values on the "Forecast Comparison" sheet, and the corresponding rolled-up
# totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------------
# Row 2 (W8)
$wsForecast.Range("D2").Value = 117
$wsForecast.Range("H2").Value = 3.8
$wsForecast.Range("L2").Value = 1.2

# Row 3 (W9)
$wsForecast.Range("D3").Value = 117
$wsForecast.Range("H3").Value = 2.8
$wsForecast.Range("L3").Value = 1.06

# Row 4 (W10)
$wsForecast.Range("H4").Value = 2.25
$wsForecast.Range("L4").Value = 0.92

# Row 5 (W11)
$wsForecast.Range("H5").Value = 1.27
$wsForecast.Range("L5").Value = 0.84

# Row 6 (W12)
$wsForecast.Range("D6").Value = 93
$wsForecast.Range("H6").Value = 0.27
$wsForecast.Range("I6").Value = "High"
$wsForecast.Range("L6").Value = 0.99

# Row 7 (W13)
$wsForecast.Range("D7").Value = 89
$wsForecast.Range("L7").Value = 1.05

# Row 8 (W14)
$wsForecast.Range("D8").Value = 93
$wsForecast.Range("L8").Value = 1.14

# Row 9 (W15)
$wsForecast.Range("D9").Value = 92
$wsForecast.Range("L9").Value = 1.19

# Row 10 (W16)
$wsForecast.Range("D10").Value = 92
$wsForecast.Range("L10").Value = 0.88

# Row 11 (W17)
$wsForecast.Range("D11").Value = 89
$wsForecast.Range("L11").Value = 1.02

# Row 12 (W18)
$wsForecast.Range("D12").Value = 91
$wsForecast.Range("L12").Value = 0.85

# Row 13 (W19)
$wsForecast.Range("D13").Value = 92
$wsForecast.Range("L13").Value = 0.86

# Row 14 (W20)
$wsForecast.Range("L14").Value = 0.9

# Row 15 (W21)
$wsForecast.Range("L15").Value = 0.82

# Row 16 (W22)
$wsForecast.Range("L16").Value = 0.89

# Row 17 (W23)
$wsForecast.Range("L17").Value = 1.01

# --- Summary sheet --------------------------------------------------------
# These cells hold their numbers as text (not numeric values), so force a
# text number format before assigning, to avoid Excel auto-converting the
# numeric-looking strings into real numbers.
$summaryTextCells = "B9", "B10", "B11", "B12", "B14"
foreach ($addr in $summaryTextCells) {
    $wsSummary.Range($addr).NumberFormat = "@"
}

$wsSummary.Range("B9").Value  = "1496"
$wsSummary.Range("B10").Value = "789"
$wsSummary.Range("B11").Value = "420"
$wsSummary.Range("B12").Value = "117"
$wsSummary.Range("B14").Value = "84"
